$wb = $excel.ActiveWorkbook

# Mapping of row -> new value for column F ("想去人数") that must be updated.
$updates = @{
    2  = 192
    3  = 3039
    4  = 219
    5  = 114
    8  = 1616
    10 = 357
    13 = 188
    15 = 224
    18 = 21
    21 = 11
    22 = 360
    23 = 161
    26 = 2032
    28 = 460
    30 = 182
    35 = 494
}

# Both the "展览" sheet and the "全部类型" sheet carry identical data and
# both need the same F-column values updated.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
